# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E) and "Correspond Handback DateTime" (H)
# timestamps for the f4bc0fdc... and f8a4f9c9... rows on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4").Value = "2016-03-18 10:16:24"
$wsZh.Range("E5").Value = "2016-03-18 10:16:24"
$wsZh.Range("H4").Value = "2016-03-18 10:16:56"
$wsZh.Range("H5").Value = "2016-03-18 10:16:56"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4").Value = "2016-03-18 10:16:27"
$wsDe.Range("E5").Value = "2016-03-18 10:16:27"
$wsDe.Range("H4").Value = "2016-03-18 10:17:02"
$wsDe.Range("H5").Value = "2016-03-18 10:17:02"
